# Apply updated crypto price/volume data per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.423.67'
$ws.Range('E2').Value = '  +3.65%  '
$ws.Range('D3').Value = '2.406.65'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.74%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.22'
$ws.Range('E5').Value = '  +1.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.40'
$ws.Range('E6').Value = '  +4.80%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.540'
$ws.Range('E8').Value = '  +2.07%  '
$ws.Range('D9').Value = '2.433.77'
$ws.Range('E9').Value = '  +2.30%  '
$ws.Range('E10').Value = '  +4.23%  '
$ws.Range('E11').Value = '  +0.66%  '
$ws.Range('E12').Value = '  +2.33%  '
$ws.Range('E13').Value = '  +3.75%  '
$ws.Range('E14').Value = '  +6.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000176'
$ws.Range('E15').Value = '  +5.38%  '
$ws.Range('D17').Value = '62.522.23'
$ws.Range('E17').Value = '  +4.00%  '
$ws.Range('D18').Value = '2.433.65'
$ws.Range('E18').Value = '  +2.39%  '
$ws.Range('E19').Value = '  -2.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.94'
$ws.Range('E20').Value = '  +3.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '326.85'
$ws.Range('E21').Value = '  +1.43%  '
$ws.Range('E22').Value = '  +1.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.997'
$ws.Range('E24').Value = '  -0.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.47'
$ws.Range('E25').Value = '  +2.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '626.65'
$ws.Range('E26').Value = '  +12.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.45'
$ws.Range('E27').Value = '  +3.71%  '
$ws.Range('D28').Value = '0.0₃0980'
$ws.Range('E28').Value = '  +5.23%  '
$ws.Range('D29').Value = '2.560.57'
$ws.Range('E29').Value = '  +2.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.16'
$ws.Range('E30').Value = '  +1.84%  '
$ws.Range('E31').Value = '  +7.22%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.85'
$ws.Range('E32').Value = '  +2.82%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.137'
$ws.Range('E33').Value = '  +3.66%  '
$ws.Range('D34').Value = '0.0₆0386'
$ws.Range('E34').Value = '  +33.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.48'
$ws.Range('E35').Value = '  +2.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.995'
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('E37').Value = '  +3.99%  '
$ws.Range('E38').Value = '  +1.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '151.96'
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.37'
$ws.Range('E40').Value = '  +6.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.62'
$ws.Range('E41').Value = '  +2.42%  '
$ws.Range('E42').Value = '  +12.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.75'
$ws.Range('E43').Value = '  +5.74%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.82'
$ws.Range('E45').Value = '  +26.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '144.62'
$ws.Range('E46').Value = '  +2.81%  '
$ws.Range('E47').Value = '  +1.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.51'
$ws.Range('E48').Value = '  +6.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.599'
$ws.Range('E49').Value = '  +1.95%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0515'
$ws.Range('E50').Value = '  +2.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0918'
$ws.Range('E51').Value = '  +2.21%  '
